$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each value is prefixed with a literal leading apostrophe (Excel's
# text-force quote prefix) so values that look numeric (e.g. "324.23",
# "0.00001114") are stored as exact text instead of being coerced to
# floating point numbers, matching the original inlineStr cell content.
$ws.Range("D2").Value = "'30.421.53"
$ws.Range("E2").Value = "'  +1.33%  "
$ws.Range("D3").Value = "'2.002.14"
$ws.Range("E3").Value = "'  +4.50%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'324.23"
$ws.Range("E5").Value = "'  +1.30%  "
$ws.Range("E6").Value = "'  +0.02%  "
$ws.Range("E7").Value = "'  +1.55%  "
$ws.Range("D8").Value = "'0.4135"
$ws.Range("E8").Value = "'  +2.61%  "
$ws.Range("D9").Value = "'0.08730"
$ws.Range("E9").Value = "'  +6.07%  "
$ws.Range("E10").Value = "'  +2.35%  "
$ws.Range("B11").Value = "'Solana"
$ws.Range("C11").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'24.54"
$ws.Range("E11").Value = "'  +3.39%  "
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.995.99"
$ws.Range("E12").Value = "'  +4.58%  "
$ws.Range("B13").Value = "'Polkadot"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.563"
$ws.Range("E13").Value = "'  +2.35%  "
$ws.Range("B14").Value = "'Chainlink"
$ws.Range("C14").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.451"
$ws.Range("B15").Value = "'BinanceUSD"
$ws.Range("C15").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "'  +0.05%  "
$ws.Range("B16").Value = "'Litecoin"
$ws.Range("C16").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'94.20"
$ws.Range("E16").Value = "'  +2.21%  "
$ws.Range("B17").Value = "'ShibaInu"
$ws.Range("C17").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001114"
$ws.Range("E17").Value = "'  +1.63%  "
$ws.Range("B18").Value = "'TRON"
$ws.Range("C18").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06509"
$ws.Range("E18").Value = "'  +0.26%  "
$ws.Range("B19").Value = "'Avalanche"
$ws.Range("C19").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'18.89"
$ws.Range("E19").Value = "'  +4.04%  "
$ws.Range("B20").Value = "'Dai"
$ws.Range("C20").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "'  +0.03%  "
$ws.Range("B21").Value = "'Uniswap"
$ws.Range("C21").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.195"
$ws.Range("E21").Value = "'  +4.37%  "
$ws.Range("B22").Value = "'WrappedBTC"
$ws.Range("C22").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "'30.476.92"
$ws.Range("E22").Value = "'  +1.38%  "
$ws.Range("B23").Value = "'Cosmos"
$ws.Range("C23").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'11.91"
$ws.Range("E23").Value = "'  +5.61%  "
$ws.Range("B24").Value = "'Toncoin"
$ws.Range("C24").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.223"
$ws.Range("E24").Value = "'  +1.36%  "
$ws.Range("B25").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "'2.230.54"
$ws.Range("E25").Value = "'  +4.70%  "
$ws.Range("B26").Value = "'EthereumClassic"
$ws.Range("C26").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'22.32"
$ws.Range("E26").Value = "'  +0.12%  "
$ws.Range("B27").Value = "'Monero"
$ws.Range("C27").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'163.02"
$ws.Range("E27").Value = "'  +0.62%  "
$ws.Range("B28").Value = "'LidoDAOToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.405"
$ws.Range("E28").Value = "'  +5.08%  "
$ws.Range("B29").Value = "'BitcoinCash"
$ws.Range("C29").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'131.11"
$ws.Range("E29").Value = "'  +1.66%  "
$ws.Range("B30").Value = "'ImmutableX"
$ws.Range("C30").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.134"
$ws.Range("E30").Value = "'  +0.37%  "
$ws.Range("B31").Value = "'Stellar"
$ws.Range("C31").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1052"
$ws.Range("E31").Value = "'  +1.32%  "
$ws.Range("B32").Value = "'Filecoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.062"
$ws.Range("E32").Value = "'  +0.93%  "
$ws.Range("B33").Value = "'HuobiToken"
$ws.Range("C33").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.829"
$ws.Range("E33").Value = "'  +0.50%  "
$ws.Range("B34").Value = "'ARBITRUM"
$ws.Range("C34").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.341"
$ws.Range("E34").Value = "'  +11.42%  "
$ws.Range("B35").Value = "'VeChain"
$ws.Range("C35").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.02520"
$ws.Range("E35").Value = "'  +3.22%  "
$ws.Range("B36").Value = "'InternetComputer(DFINITY)"
$ws.Range("C36").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.431"
$ws.Range("E36").Value = "'  +1.63%  "
$ws.Range("B37").Value = "'Hedera"
$ws.Range("C37").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06596"
$ws.Range("E37").Value = "'  +2.61%  "
$ws.Range("B38").Value = "'Aptos"
$ws.Range("C38").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'12.48"
$ws.Range("E38").Value = "'  +9.65%  "
$ws.Range("B39").Value = "'Algorand"
$ws.Range("C39").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2196"
$ws.Range("E39").Value = "'  +1.56%  "
$ws.Range("B40").Value = "'FraxShare"
$ws.Range("C40").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'9.026"
$ws.Range("E40").Value = "'  +1.38%  "
$ws.Range("B41").Value = "'TheSandbox"
$ws.Range("C41").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6628"
$ws.Range("E41").Value = "'  +2.93%  "
$ws.Range("B42").Value = "'TrustWalletToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.230"
$ws.Range("E42").Value = "'  +1.05%  "
$ws.Range("B43").Value = "'EnergySwap"
$ws.Range("C43").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'13.56"
$ws.Range("E43").Value = "'  +1.08%  "
$ws.Range("B44").Value = "'Decentraland"
$ws.Range("C44").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6161"
$ws.Range("E44").Value = "'  +2.76%  "
$ws.Range("B45").Value = "'NEARProtocol"
$ws.Range("C45").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.190"
$ws.Range("E45").Value = "'  -0.23%  "
$ws.Range("B46").Value = "'PancakeSwap"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.666"
$ws.Range("E46").Value = "'  +0.87%  "
$ws.Range("B47").Value = "'EOS"
$ws.Range("C47").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D47").Value = "'1.261"
$ws.Range("E47").Value = "'  +3.91%  "
$ws.Range("B48").Value = "'Quant"
$ws.Range("C48").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'124.35"
$ws.Range("E48").Value = "'  +0.75%  "
$ws.Range("B49").Value = "'Aave"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'80.24"
$ws.Range("E49").Value = "'  +1.87%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06887"
$ws.Range("E50").Value = "'  +1.35%  "
$ws.Range("B51").Value = "'WEMIXTOKEN"
$ws.Range("C51").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'1.105"
$ws.Range("E51").Value = "'  -2.15%  "
